$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 627.9091
$ws.Range("I8").Value = 34.166668
$ws.Range("K8").Value = 102.500004
$ws.Range("M8").Value = 36.499996
# Row 34
$ws.Range("H34").Value = 6665
$ws.Range("I34").Value = 6665
$ws.Range("K34").Value = 6665
$ws.Range("M34").Value = -6462
# Row 36
$ws.Range("H36").Value = 6665
$ws.Range("I36").Value = 6665
$ws.Range("K36").Value = 6665
$ws.Range("M36").Value = -5950

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 50
$ws.Range("H50").Value = 733
$ws.Range("I50").Value = 200
$ws.Range("J50").Value = 999.5
$ws.Range("K50").Value = 200
$ws.Range("L50").Value = 999.5
$ws.Range("M50").Value = 514
$ws.Range("N50").Value = -2427.5
# Row 74
$ws.Range("H74").Value = 2165.625
$ws.Range("I74").Value = 2351.3157
$ws.Range("J74").Value = 1460
$ws.Range("K74").Value = 2351.3157
$ws.Range("L74").Value = 1460
$ws.Range("M74").Value = -1477.3157
$ws.Range("N74").Value = -3208
# Row 77
$ws.Range("H77").Value = 2165.625
$ws.Range("I77").Value = 2351.3157
$ws.Range("J77").Value = 1460
$ws.Range("K77").Value = 11756.5785
$ws.Range("L77").Value = 7300
$ws.Range("M77").Value = -7388.5785
$ws.Range("N77").Value = -16036
# Row 122
$ws.Range("H122").Value = 4127.75
$ws.Range("I122").Value = 3837
$ws.Range("K122").Value = 11511
$ws.Range("M122").Value = -9061

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6874.615
$ws.Range("I86").Value = 4624.7144
$ws.Range("K86").Value = 4624.7144
$ws.Range("M86").Value = -3501.7144
# Row 89
$ws.Range("H89").Value = 6874.615
$ws.Range("I89").Value = 4624.7144
$ws.Range("K89").Value = 23123.572
$ws.Range("M89").Value = -17507.572

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6635.893
$ws.Range("I31").Value = 996.3333
$ws.Range("J31").Value = 7312.64
$ws.Range("K31").Value = 996.3333
$ws.Range("L31").Value = 7312.64
$ws.Range("M31").Value = -701.3333
$ws.Range("N31").Value = -7902.64
# Row 34
$ws.Range("H34").Value = 6635.893
$ws.Range("I34").Value = 996.3333
$ws.Range("J34").Value = 7312.64
$ws.Range("K34").Value = 996.3333
$ws.Range("L34").Value = 7312.64
$ws.Range("M34").Value = -794.3333
$ws.Range("N34").Value = -7716.64
# Row 70
$ws.Range("H70").Value = 40090
$ws.Range("J70").Value = 40090
$ws.Range("L70").Value = 40090
$ws.Range("N70").Value = -40720
# Row 73
$ws.Range("H73").Value = 40090
$ws.Range("J73").Value = 40090
$ws.Range("L73").Value = 40090
$ws.Range("N73").Value = -42274
# Row 80
$ws.Range("H80").Value = 55128
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 55128
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
# Row 107
$ws.Range("H107").Value = 630
$ws.Range("I107").Value = 575
$ws.Range("K107").Value = 575
$ws.Range("M107").Value = 1345

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 449.92856
$ws.Range("I7").Value = 407
$ws.Range("K7").Value = 1221
$ws.Range("M7").Value = -1109
# Row 17
$ws.Range("H17").Value = 3804.1538
$ws.Range("J17").Value = 9847.6
$ws.Range("L17").Value = 29542.8
$ws.Range("N17").Value = -29880.8
# Row 18
$ws.Range("H18").Value = 2854.3333
$ws.Range("I18").Value = 397
$ws.Range("K18").Value = 1191
$ws.Range("M18").Value = -1022
# Row 34
$ws.Range("H34").Value = 899.1818
$ws.Range("I34").Value = 214
$ws.Range("K34").Value = 642
$ws.Range("M34").Value = -558
# Row 39
$ws.Range("H39").Value = 6383.5713
$ws.Range("J39").Value = 7395
$ws.Range("L39").Value = 22185
$ws.Range("N39").Value = -22773
# Row 41
$ws.Range("H41").Value = 3383.6667
$ws.Range("I41").Value = 75
$ws.Range("J41").Value = 10001
$ws.Range("K41").Value = 225
$ws.Range("L41").Value = 30003
$ws.Range("M41").Value = 113
$ws.Range("N41").Value = -30679
# Row 55
$ws.Range("H55").Value = 8538.833000000001
$ws.Range("J55").Value = 10335.75
$ws.Range("L55").Value = 31007.25
$ws.Range("N55").Value = -31361.25
# Row 70
$ws.Range("H70").Value = 1400
$ws.Range("I70").Value = 1400
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4200
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3885
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 1400
$ws.Range("I73").Value = 1400
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4200
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3108
$ws.Range("N73").ClearContents()
# Row 131
$ws.Range("H131").Value = 1580.3478
$ws.Range("I131").Value = 893.8182
$ws.Range("K131").Value = 2681.4546
$ws.Range("M131").Value = 2358.5454
# Row 132
$ws.Range("H132").Value = 1558.125
$ws.Range("J132").Value = 1446
$ws.Range("L132").Value = 13014
$ws.Range("N132").Value = -18074

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 12499
$ws.Range("I80").Value = 9998
$ws.Range("J80").Value = 15000
$ws.Range("K80").Value = 9998
$ws.Range("L80").Value = 15000
$ws.Range("M80").Value = -9000
$ws.Range("N80").Value = -16996
# Row 83
$ws.Range("H83").Value = 12499
$ws.Range("I83").Value = 9998
$ws.Range("J83").Value = 15000
$ws.Range("K83").Value = 49990
$ws.Range("L83").Value = 75000
$ws.Range("M83").Value = -44998
$ws.Range("N83").Value = -84984
# Row 126
$ws.Range("H126").Value = 4000
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1025
$ws.Range("I93").Value = 875
$ws.Range("K93").Value = 875
$ws.Range("M93").Value = 373

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 50
$ws.Range("H50").Value = 14996
$ws.Range("J50").Value = 14996
$ws.Range("L50").Value = 14996
$ws.Range("N50").Value = -16258
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 126
$ws.Range("H126").Value = 5223
$ws.Range("I126").Value = 3515.25
$ws.Range("J126").Value = 6930.75
$ws.Range("K126").Value = 10545.75
$ws.Range("L126").Value = 20792.25
$ws.Range("M126").Value = -8075.75
$ws.Range("N126").Value = -25732.25
